# DOCS: applied all pending changes.
$wb = $excel.ActiveWorkbook

# --- Workbook / sheet-level changes -----------------------------------
# Rename the two existing sheets and append a new blank "Sheet1" at the end.
$wsOld = $wb.Worksheets.Item(1)
$wsNew = $wb.Worksheets.Item(2)
$wsOld.Name = "Red Math Ops (old)"
$wsNew.Name = "Red Math Ops"

$wsBlank = $wb.Worksheets.Add()
$wsBlank.Name = "Sheet1"
$wsBlank.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

Write-Host "done"
